# Update odds data on the active worksheet (row 3 and row 5)
# as described in the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates ---
$ws.Range("G3").Value  = 2
$ws.Range("H3").Value  = 3.4
$ws.Range("I3").Value  = 3.3
$ws.Range("L3").Value  = 4
$ws.Range("U3").Value  = 1.73
$ws.Range("V3").Value  = 2
$ws.Range("W3").Value  = 8
$ws.Range("X3").Value  = 10
$ws.Range("Z3").Value  = 19
$ws.Range("AG3").Value = 201
$ws.Range("AI3").Value = 19
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 26
$ws.Range("AM3").Value = 34
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 11
$ws.Range("AW3").Value = 5.5
$ws.Range("AX3").Value = 19
$ws.Range("AY3").Value = 26
$ws.Range("BA3").Value = 81

# --- Row 5 updates ---
$ws.Range("P5").Value = 3.6
$ws.Range("U5").Value = 1.6
$ws.Range("V5").Value = 2.27
